$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Insert a new product row at row 5 ("KERELLA LOTION 30 ML"), between
#    the existing "ELICA-M CREAM 30 GRAM" (row 4) and
#    "LIDOCAINE 10% TOPICAL SPRAY 15 GM" (originally row 5, now row 6).
# ---------------------------------------------------------------------
$ws.Rows.Item(5).Insert()
$ws.Range("A6:N6").Copy()
$ws.Range("A5:N5").PasteSpecial($xlPasteFormats)
$ws.Rows.Item(5).RowHeight = 25.5
$ws.Range("B5:G5").Merge()
$ws.Range("H5:K5").Merge()
$ws.Range("L5:M5").Merge()

$ws.Range("B5").Value = "KERELLA LOTION 30 ML"
$ws.Range("H5").Value = "3:0"
$ws.Range("L5").Value = 31
$ws.Range("N5").Value = "1:0"

# ---------------------------------------------------------------------
# 2) Insert a new product row at row 8 ("SUPOLACK HAIR SHAMPOO 200 ML"),
#    between "NEUROGLOPENTIN 300 MG 30 CAPS." (now row 7) and
#    "TELFAST 180MG 20 F.C. TABS" (currently row 8, about to shift to 9).
# ---------------------------------------------------------------------
$ws.Rows.Item(8).Insert()
$ws.Range("A7:N7").Copy()
$ws.Range("A8:N8").PasteSpecial($xlPasteFormats)
$ws.Rows.Item(8).RowHeight = 25.5
$ws.Range("B8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()

$ws.Range("B8").Value = "SUPOLACK HAIR SHAMPOO 200 ML"
$ws.Range("H8").Value = "0:0"
$ws.Range("L8").Value = 149.5
$ws.Range("N8").Value = "1:0"

# ---------------------------------------------------------------------
# Re-sequence the "م" (serial number) column for the full product table
# (rows 4..15) and re-write the totals / footer rows that shifted down
# from 14/15 to 16/17.
# ---------------------------------------------------------------------
for ($i = 0; $i -lt 12; $i++) {
    $r = 4 + $i
    $ws.Cells.Item($r, 1).Value = $i + 1
}

# Grand total of the "سعر البيع" column (L) across all 12 product rows.
$ws.Range("K16").Value = 657.5
